$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 (Grid) : update the "Http" error comment text ---
$ws.Range("I12").Value = "KO  : System.InvalidOperationException: Cannot provide a value for property 'Http'"

# --- Row 14 (Dialog) : add comment, switch to red/wrap style, taller row ---
$ws.Range("I14").Value = "Mitigé : fonctionne très bien mais nécessite d'ajouter le code de boite de dialogue dans la page HTML"
$ws.Range("I14").Font.Color = 255
$ws.Range("I14").Font.Bold = $true
$ws.Range("I14").WrapText = $true
$ws.Rows("14").RowHeight = 30

# --- Row 15 (Toast) : simple OK ---
$ws.Range("I15").Value = "OK"

# --- Row 16 (Chart) : simple NA ---
$ws.Range("I16").Value = "NA"

# --- Row 17 (Scheduler) : simple NA ---
$ws.Range("I17").Value = "NA"

# --- Row 18 (Autocomplete) : add comment, switch to red/wrap style, taller row ---
$ws.Range("I18").Value = "KO : Il existe un autocomplete mais je ne sais pas trop s'il sait gérer le chargement dynamique. Tous les exemples sont permettent uniquement de rechercher les valeurs en saisissant un texte."
$ws.Range("I18").Font.Color = 255
$ws.Range("I18").Font.Bold = $true
$ws.Range("I18").WrapText = $true
$ws.Rows("18").RowHeight = 45

# --- Row 19 (ListBox) : add comment, switch to red/wrap style, taller row ---
$ws.Range("I19").Value = "Binding : KO (Value de type string uniquement + les données s'affichent mais je n'arrive pas à binder la valeur à ma propriété)"
$ws.Range("I19").Font.Color = 255
$ws.Range("I19").Font.Bold = $true
$ws.Range("I19").WrapText = $true
$ws.Rows("19").RowHeight = 30

# --- Update the active selection to match where the author ended up working ---
$ws.Range("I18").Select()
